$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($null -ne $val -and $val -ne "") {
        $parts = @($val -split ",\s*")
        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $hasSystem = $true }
        }
        if ($hasSystem) {
            $rest = @()
            $removed = $false
            foreach ($p in $parts) {
                if ((-not $removed) -and $p.Equals("System")) {
                    $removed = $true
                    continue
                }
                $rest += $p
            }
            $newParts = @("System") + $rest
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value = $newVal
        }
    }
}
